$d = $word.ActiveDocument

# --- Change 1a: first placeholder "${cscAdviser}" -> "${cscPresident}" ---
$rng = $d.Content
$rng.Find.Execute('${cscAdviser}', $false, $false, $false, $false, $false, $true, 1, $false, '${cscPresident}', 2)

# --- Change 1b: the 45-space run between the tabs and "${oicOsa}" shrinks to 13 spaces ---
$rng = $d.Content
$rng.Find.Execute('                                             ', $false, $false, $false, $false, $false, $true, 1, $false, '             ', 2)

# --- Change 1c: second placeholder "${oicOsa}" -> "${cscAdviser}" ---
$rng = $d.Content
$rng.Find.Execute('${oicOsa}', $false, $false, $false, $false, $false, $true, 1, $false, '${cscAdviser}', 2)

# --- Change 2a: "                    Adviser" (20 spaces + Adviser) -> " Adviser" (1 space + Adviser) ---
$rng = $d.Content
$rng.Find.Execute('                    Adviser', $false, $false, $false, $false, $false, $true, 1, $false, ' Adviser', 2)

# --- Change 2b: insert 5 additional tabs right after " Adviser" (before the existing trailing tabs) ---
$rng = $d.Content
$rng.Find.Execute(' Adviser', $false, $false, $false, $false, $false, $true, 1, $false, '', 0)
$insertPos = $rng.End
$insertRange = $d.Range($insertPos, $insertPos)
$insertRange.InsertAfter("`t`t`t`t`t")

Write-Output "done"
